# Streamlit Approval System - append two new pending-approval rows
# (WGG 02 / Western Interior Designers & Marine Contractors) to Sheet1,
# as rows 3 and 4 below the existing header + single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 --------------------------------------------------------------
$ws.Range("A3").Value = "WGG 02"
$ws.Range("B3").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C3").Value = 46297
$ws.Range("C3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D3").Value = 286962
$ws.Range("E3").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F3").Value = 34413429360
$ws.Range("G3").Value = "NEFT"
$ws.Range("H3").Value = "SBIN0003229"
$ws.Range("I3").Value = "AAAFW8862C"
$ws.Range("J3").Value = "32AAAFW8862C1Z9"
$ws.Range("L3").Value = "eed77dce-c7f0-4070-9404-e16ea6c44ea4"
$ws.Range("U3").Value = "pending"
$ws.Range("V3").Value = 45560
$ws.Range("X3").Value = "2 month consultancy fees and Hisham sir & Hijas Sir tax tax repayment (45000+560) RPA_UNIQUE_ID : 241b1377-30af-494d-90b7-80778b59d820"
$ws.Range("Y3").Value = "Cochin"
$ws.Range("Z3").Value = "PAYMENT"
$ws.Range("AA3").Value = "Payments@westernidc.com"
$ws.Range("AB3").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC3").Value = 0
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0

# ---- Row 4 --------------------------------------------------------------
$ws.Range("A4").Value = "WGG 02"
$ws.Range("B4").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C4").Value = 46297
$ws.Range("C4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D4").Value = 286962
$ws.Range("E4").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F4").Value = 34413429360
$ws.Range("G4").Value = "NEFT"
$ws.Range("H4").Value = "SBIN0003229"
$ws.Range("I4").Value = "AAAFW8862C"
$ws.Range("J4").Value = "32AAAFW8862C1Z9"
$ws.Range("L4").Value = "4e75250f-bee6-47ed-8e72-941a23dfdecc"
$ws.Range("U4").Value = "pending"
$ws.Range("V4").Value = 40000
$ws.Range("X4").Value = "Shabeena Beevi RPA_UNIQUE_ID : 45244746-a8e8-45ec-8285-f14adba82cac"
$ws.Range("Y4").Value = "Cochin"
$ws.Range("Z4").Value = "PAYMENT"
$ws.Range("AA4").Value = "Payments@westernidc.com"
$ws.Range("AB4").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC4").Value = 0
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
